$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Cason Wallace"
$ws.Cells.Item(2, 2).Value = "PG,SG"
$ws.Cells.Item(2, 3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(3, 1).Value = "Jordan Poole"
$ws.Cells.Item(3, 2).Value = "PG,SG"
$ws.Cells.Item(3, 3).Value = "Washington Wizards"
$ws.Cells.Item(4, 1).Value = "Jamal Murray"
$ws.Cells.Item(4, 2).Value = "PG,SG"
$ws.Cells.Item(4, 3).Value = "Denver Nuggets"
$ws.Cells.Item(5, 1).Value = "Shai Gilgeous-Alexander"
$ws.Cells.Item(5, 2).Value = "PG,SG"
$ws.Cells.Item(5, 3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(6, 1).Value = "Naji Marshall"
$ws.Cells.Item(6, 2).Value = "SG,SF,PF"
$ws.Cells.Item(6, 3).Value = "Dallas Mavericks"
$ws.Cells.Item(7, 1).Value = "Dillon Brooks"
$ws.Cells.Item(7, 2).Value = "SG,SF,PF"
$ws.Cells.Item(7, 3).Value = "Houston Rockets"
$ws.Cells.Item(8, 1).Value = "Zach LaVine"
$ws.Cells.Item(8, 2).Value = "SG,SF"
$ws.Cells.Item(8, 3).Value = "Sacramento Kings"
$ws.Cells.Item(9, 1).Value = "Buddy Hield"
$ws.Cells.Item(9, 2).Value = "SG,SF"
$ws.Cells.Item(9, 3).Value = "Golden State Warriors"
$ws.Cells.Item(10, 1).Value = "Donovan Clingan"
$ws.Cells.Item(10, 2).Value = "C"
$ws.Cells.Item(10, 3).Value = "Portland Trail Blazers"
$ws.Cells.Item(11, 1).Value = "Alperen Sengün"
$ws.Cells.Item(11, 2).Value = "C"
$ws.Cells.Item(11, 3).Value = "Houston Rockets"
$ws.Cells.Item(12, 1).Value = "Kyle Kuzma"
$ws.Cells.Item(12, 2).Value = "SF,PF"
$ws.Cells.Item(12, 3).Value = "Milwaukee Bucks"
$ws.Cells.Item(13, 1).Value = "Kris Dunn"
$ws.Cells.Item(13, 2).Value = "PG,SG"
$ws.Cells.Item(13, 3).Value = "LA Clippers"
$ws.Cells.Item(14, 1).Value = "Khris Middleton"
$ws.Cells.Item(14, 2).Value = "SF"
$ws.Cells.Item(14, 3).Value = "Washington Wizards"
$ws.Cells.Item(15, 1).Value = "Dennis Schröder"
$ws.Cells.Item(15, 2).Value = "PG,SG"
$ws.Cells.Item(15, 3).Value = "Detroit Pistons"
$ws.Cells.Item(16, 1).Value = "Rui Hachimura"
$ws.Cells.Item(16, 2).Value = "SF,PF"
$ws.Cells.Item(16, 3).Value = "Los Angeles Lakers"
$ws.Cells.Item(17, 1).Value = "Jordan Clarkson"
$ws.Cells.Item(17, 2).Value = "SG,SF"
$ws.Cells.Item(17, 3).Value = "Utah Jazz"
$ws.Cells.Item(18, 1).Value = "John Collins"
$ws.Cells.Item(18, 2).Value = "PF,C"
$ws.Cells.Item(18, 3).Value = "Utah Jazz"
